# Update column F ("dSF") values on the active sheet to match the
# re-pulled / re-calculated data described in the commit
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = 1
    5  = -3
    7  = 1
    8  = -1
    9  = -3
    11 = 0
    14 = 2
    16 = 3
    20 = -1
    21 = -9
    22 = -10
    26 = 3
    28 = -5
    38 = 1
    41 = 2
    43 = -2
    47 = 2
    48 = 5
    52 = 6
    54 = 6
    56 = -5
    59 = 2
    60 = 4
    61 = -1
    65 = -1
    70 = 2
    71 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
